# Update "想去人数" (F column) values across the four sheets to reflect
# newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 1158
$ws.Cells.Item(6, 6).Value = 2776
$ws.Cells.Item(7, 6).Value = 232
$ws.Cells.Item(9, 6).Value = 104
$ws.Cells.Item(10, 6).Value = 291
$ws.Cells.Item(13, 6).Value = 110
$ws.Cells.Item(14, 6).Value = 135
$ws.Cells.Item(15, 6).Value = 1744
$ws.Cells.Item(18, 6).Value = 206
$ws.Cells.Item(19, 6).Value = 259

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 18
$ws.Cells.Item(10, 6).Value = 39
$ws.Cells.Item(23, 6).Value = 29

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 6362
$ws.Cells.Item(4, 6).Value = 2032
$ws.Cells.Item(5, 6).Value = 271

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 6362
$ws.Cells.Item(4, 6).Value = 2032
$ws.Cells.Item(5, 6).Value = 271
$ws.Cells.Item(12, 6).Value = 1158
$ws.Cells.Item(14, 6).Value = 18
$ws.Cells.Item(17, 6).Value = 2776
$ws.Cells.Item(19, 6).Value = 232
$ws.Cells.Item(20, 6).Value = 39
$ws.Cells.Item(24, 6).Value = 104
$ws.Cells.Item(25, 6).Value = 291
$ws.Cells.Item(29, 6).Value = 110
$ws.Cells.Item(30, 6).Value = 135
$ws.Cells.Item(32, 6).Value = 1744
$ws.Cells.Item(37, 6).Value = 206
$ws.Cells.Item(44, 6).Value = 29
$ws.Cells.Item(45, 6).Value = 259

$wb.Save()
